$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run the "Sort" the author performed on the data range A2:T20,
# ascending by column J (p_val) - this reorders rows 2-20.
$sortRange = $ws.Range("A2:T20")
$sortKey = $ws.Range("J2:J20")
$sortRange.Sort($sortKey, 1)

# Update the AutoFilter: domain (column A / colId 0) filter changes
# from "biological" to "chemical". Re-apply the model (column C / colId 2)
# filter afterwards so it keeps its original relative ordering in the XML.
$ws.Range("A1:T58").AutoFilter(1, @("chemical"), 7)
$ws.Range("A1:T58").AutoFilter(3, @("field_in_loc_as_random"), 7)
